$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column CK: experience_group
$ws.Range("CK1").Value = "experience_group"
$ws.Range("CJ1").Copy() | Out-Null
$ws.Range("CK1").PasteSpecial(-4122) | Out-Null

# Per-row experience_group values
$ws.Range("CK2").Value = "Novice"
$ws.Range("CK3").Value = "Novice"
$ws.Range("CK4").Value = "Novice"
$ws.Range("CK5").Value = "Novice+"
$ws.Range("CK6").Value = "Novice+"
$ws.Range("CK7").Value = "Novice"
$ws.Range("CK8").Value = "Novice+"
$ws.Range("CK9").Value = "Novice+"
$ws.Range("CK10").Value = "Novice"
$ws.Range("CK11").Value = "Novice"
$ws.Range("CK12").Value = "Novice"
$ws.Range("CK13").Value = "Network/IT admin"
$ws.Range("CK14").Value = "Novice"
$ws.Range("CK15").Value = "Novice"
$ws.Range("CK16").Value = "Novice+"
$ws.Range("CK17").Value = "Novice"
$ws.Range("CK18").Value = "Novice"
$ws.Range("CK19").Value = "Novice"
$ws.Range("CK20").Value = "Novice"
$ws.Range("CK21").Value = "Novice"
$ws.Range("CK22").Value = "Novice"
$ws.Range("CK23").Value = "Novice"
$ws.Range("CK24").Value = "Novice"
$ws.Range("CK25").Value = "Novice"
$ws.Range("CK26").Value = "Novice"
$ws.Range("CK27").Value = "Novice"
$ws.Range("CK28").Value = "Novice"
$ws.Range("CK29").Value = "Novice"
$ws.Range("CK30").Value = "Novice"
$ws.Range("CK31").Value = "Novice"
$ws.Range("CK32").Value = "Novice"
$ws.Range("CK33").Value = "Novice"
$ws.Range("CK34").Value = "Novice"
$ws.Range("CK35").Value = "Novice"
$ws.Range("CK36").Value = "Cyber security"
$ws.Range("CK37").Value = "Novice"
$ws.Range("CK38").Value = "Novice"
$ws.Range("CK39").Value = "Novice"
$ws.Range("CK40").Value = "Novice"
$ws.Range("CK41").Value = "Novice"
$ws.Range("CK42").Value = "Novice"
$ws.Range("CK43").Value = "Novice"
$ws.Range("CK44").Value = "Novice"
$ws.Range("CK45").Value = "Novice"
$ws.Range("CK46").Value = "Novice+"
$ws.Range("CK47").Value = "Novice"
$ws.Range("CK48").Value = "Cyber security"
$ws.Range("CK49").Value = "Novice"
$ws.Range("CK50").Value = "Novice"
$ws.Range("CK51").Value = "Novice"
$ws.Range("CK52").Value = "Novice+"
